$d = $word.ActiveDocument

$d.Content.Find.Execute("54+37=91", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=92", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "67-16=51", 2) | Out-Null
$d.Content.Find.Execute("50-10=40", $true, $false, $false, $false, $false, $true, 1, $false, "77-17=60", 2) | Out-Null
$d.Content.Find.Execute("77-71=6", $true, $false, $false, $false, $false, $true, 1, $false, "87+1=88", 2) | Out-Null
$d.Content.Find.Execute("1+86=87", $true, $false, $false, $false, $false, $true, 1, $false, "38+6=44", 2) | Out-Null
$d.Content.Find.Execute("65+21=86", $true, $false, $false, $false, $false, $true, 1, $false, "25+26=51", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=25", 2) | Out-Null
$d.Content.Find.Execute("26+16=42", $true, $false, $false, $false, $false, $true, 1, $false, "18+13=31", 2) | Out-Null
$d.Content.Find.Execute("3+68=71", $true, $false, $false, $false, $false, $true, 1, $false, "58-39=19", 2) | Out-Null
$d.Content.Find.Execute("38-25=13", $true, $false, $false, $false, $false, $true, 1, $false, "76+10=86", 2) | Out-Null
$d.Content.Find.Execute("23+6=29", $true, $false, $false, $false, $false, $true, 1, $false, "84+9=93", 2) | Out-Null
$d.Content.Find.Execute("0+78=78", $true, $false, $false, $false, $false, $true, 1, $false, "97-38=59", 2) | Out-Null
$d.Content.Find.Execute("15+0=15", $true, $false, $false, $false, $false, $true, 1, $false, "31+16=47", 2) | Out-Null
$d.Content.Find.Execute("61-26=35", $true, $false, $false, $false, $false, $true, 1, $false, "95-94=1", 2) | Out-Null
$d.Content.Find.Execute("34+49=83", $true, $false, $false, $false, $false, $true, 1, $false, "22+27=49", 2) | Out-Null
$d.Content.Find.Execute("6+79=85", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=77", 2) | Out-Null
$d.Content.Find.Execute("61-27=34", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=28", 2) | Out-Null
$d.Content.Find.Execute("78-34=44", $true, $false, $false, $false, $false, $true, 1, $false, "93-81=12", 2) | Out-Null
$d.Content.Find.Execute("71-34=37", $true, $false, $false, $false, $false, $true, 1, $false, "7+78=85", 2) | Out-Null
$d.Content.Find.Execute("28-21=7", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=83", 2) | Out-Null
$d.Content.Find.Execute("12+31=43", $true, $false, $false, $false, $false, $true, 1, $false, "10+43=53", 2) | Out-Null
$d.Content.Find.Execute("27+10=37", $true, $false, $false, $false, $false, $true, 1, $false, "51-22=29", 2) | Out-Null
$d.Content.Find.Execute("0+90=90", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=85", 2) | Out-Null
$d.Content.Find.Execute("53-22=31", $true, $false, $false, $false, $false, $true, 1, $false, "97-78=19", 2) | Out-Null
$d.Content.Find.Execute("39+12=51", $true, $false, $false, $false, $false, $true, 1, $false, "33+54=87", 2) | Out-Null
$d.Content.Find.Execute("81-72=9", $true, $false, $false, $false, $false, $true, 1, $false, "7+45=52", 2) | Out-Null
$d.Content.Find.Execute("3+64=67", $true, $false, $false, $false, $false, $true, 1, $false, "47-24=23", 2) | Out-Null
$d.Content.Find.Execute("13+19=32", $true, $false, $false, $false, $false, $true, 1, $false, "5+23=28", 2) | Out-Null
$d.Content.Find.Execute("63-62=1", $true, $false, $false, $false, $false, $true, 1, $false, "29+0=29", 2) | Out-Null
$d.Content.Find.Execute("40+42=82", $true, $false, $false, $false, $false, $true, 1, $false, "97-77=20", 2) | Out-Null
$d.Content.Find.Execute("78+21=99", $true, $false, $false, $false, $false, $true, 1, $false, "62+35=97", 2) | Out-Null
$d.Content.Find.Execute("64-31=33", $true, $false, $false, $false, $false, $true, 1, $false, "63+18=81", 2) | Out-Null
$d.Content.Find.Execute("18+50=68", $true, $false, $false, $false, $false, $true, 1, $false, "6+46=52", 2) | Out-Null
$d.Content.Find.Execute("92-25=67", $true, $false, $false, $false, $false, $true, 1, $false, "8+91=99", 2) | Out-Null
$d.Content.Find.Execute("75-69=6", $true, $false, $false, $false, $false, $true, 1, $false, "43-34=9", 2) | Out-Null
$d.Content.Find.Execute("68+12=80", $true, $false, $false, $false, $false, $true, 1, $false, "55+17=72", 2) | Out-Null
$d.Content.Find.Execute("46-13=33", $true, $false, $false, $false, $false, $true, 1, $false, "43-12=31", 2) | Out-Null
$d.Content.Find.Execute("38+18=56", $true, $false, $false, $false, $false, $true, 1, $false, "9-1=8", 2) | Out-Null
$d.Content.Find.Execute("89-78=11", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 2) | Out-Null
$d.Content.Find.Execute("14+18=32", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=27", 2) | Out-Null
$d.Content.Find.Execute("27+31=58", $true, $false, $false, $false, $false, $true, 1, $false, "66+14=80", 2) | Out-Null
$d.Content.Find.Execute("9+0=9", $true, $false, $false, $false, $false, $true, 1, $false, "30-7=23", 2) | Out-Null
$d.Content.Find.Execute("22+24=46", $true, $false, $false, $false, $false, $true, 1, $false, "31+5=36", 2) | Out-Null
$d.Content.Find.Execute("8+39=47", $true, $false, $false, $false, $false, $true, 1, $false, "42+40=82", 2) | Out-Null
$d.Content.Find.Execute("5+21=26", $true, $false, $false, $false, $false, $true, 1, $false, "14+14=28", 2) | Out-Null
$d.Content.Find.Execute("40-39=1", $true, $false, $false, $false, $false, $true, 1, $false, "82-4=78", 2) | Out-Null
$d.Content.Find.Execute("15+30=45", $true, $false, $false, $false, $false, $true, 1, $false, "29+61=90", 2) | Out-Null
$d.Content.Find.Execute("56+17=73", $true, $false, $false, $false, $false, $true, 1, $false, "92-83=9", 2) | Out-Null
$d.Content.Find.Execute("33-16=17", $true, $false, $false, $false, $false, $true, 1, $false, "53+32=85", 2) | Out-Null
$d.Content.Find.Execute("60+18=78", $true, $false, $false, $false, $false, $true, 1, $false, "65+28=93", 2) | Out-Null
$d.Content.Find.Execute("10+80=90", $true, $false, $false, $false, $false, $true, 1, $false, "99-84=15", 2) | Out-Null
$d.Content.Find.Execute("79-22=57", $true, $false, $false, $false, $false, $true, 1, $false, "27-15=12", 2) | Out-Null
$d.Content.Find.Execute("17+4=21", $true, $false, $false, $false, $false, $true, 1, $false, "91-21=70", 2) | Out-Null
$d.Content.Find.Execute("76-30=46", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=34", 2) | Out-Null
$d.Content.Find.Execute("54+17=71", $true, $false, $false, $false, $false, $true, 1, $false, "61+7=68", 2) | Out-Null
$d.Content.Find.Execute("80-48=32", $true, $false, $false, $false, $false, $true, 1, $false, "91-67=24", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $false, $false, $false, $false, $true, 1, $false, "74+0=74", 2) | Out-Null
$d.Content.Find.Execute("56-37=19", $true, $false, $false, $false, $false, $true, 1, $false, "17+75=92", 2) | Out-Null
$d.Content.Find.Execute("83-26=57", $true, $false, $false, $false, $false, $true, 1, $false, "28+22=50", 2) | Out-Null
$d.Content.Find.Execute("10+10=20", $true, $false, $false, $false, $false, $true, 1, $false, "95-78=17", 2) | Out-Null
$d.Content.Find.Execute("76-35=41", $true, $false, $false, $false, $false, $true, 1, $false, "81-36=45", 2) | Out-Null
$d.Content.Find.Execute("11+62=73", $true, $false, $false, $false, $false, $true, 1, $false, "37+12=49", 2) | Out-Null
$d.Content.Find.Execute("35+21=56", $true, $false, $false, $false, $false, $true, 1, $false, "45+13=58", 2) | Out-Null
$d.Content.Find.Execute("1+18=19", $true, $false, $false, $false, $false, $true, 1, $false, "18+55=73", 2) | Out-Null
$d.Content.Find.Execute("44+21=65", $true, $false, $false, $false, $false, $true, 1, $false, "18+56=74", 2) | Out-Null
$d.Content.Find.Execute("43+50=93", $true, $false, $false, $false, $false, $true, 1, $false, "87-43=44", 2) | Out-Null
$d.Content.Find.Execute("24+5=29", $true, $false, $false, $false, $false, $true, 1, $false, "83-66=17", 2) | Out-Null
$d.Content.Find.Execute("49-26=23", $true, $false, $false, $false, $false, $true, 1, $false, "23-11=12", 2) | Out-Null
$d.Content.Find.Execute("61-50=11", $true, $false, $false, $false, $false, $true, 1, $false, "72+1=73", 2) | Out-Null
$d.Content.Find.Execute("25+52=77", $true, $false, $false, $false, $false, $true, 1, $false, "83-43=40", 2) | Out-Null
$d.Content.Find.Execute("92-6=86", $true, $false, $false, $false, $false, $true, 1, $false, "40-35=5", 2) | Out-Null
$d.Content.Find.Execute("55+38=93", $true, $false, $false, $false, $false, $true, 1, $false, "53-33=20", 2) | Out-Null
$d.Content.Find.Execute("82-17=65", $true, $false, $false, $false, $false, $true, 1, $false, "80-1=79", 2) | Out-Null
$d.Content.Find.Execute("88-67=21", $true, $false, $false, $false, $false, $true, 1, $false, "88-19=69", 2) | Out-Null
$d.Content.Find.Execute("58-48=10", $true, $false, $false, $false, $false, $true, 1, $false, "27+65=92", 2) | Out-Null
$d.Content.Find.Execute("22+61=83", $true, $false, $false, $false, $false, $true, 1, $false, "96-34=62", 2) | Out-Null
$d.Content.Find.Execute("64-0=64", $true, $false, $false, $false, $false, $true, 1, $false, "1-1=0", 2) | Out-Null
$d.Content.Find.Execute("82-14=68", $true, $false, $false, $false, $false, $true, 1, $false, "88-77=11", 2) | Out-Null
$d.Content.Find.Execute("89-20=69", $true, $false, $false, $false, $false, $true, 1, $false, "38-2=36", 2) | Out-Null
$d.Content.Find.Execute("59+23=82", $true, $false, $false, $false, $false, $true, 1, $false, "45+39=84", 2) | Out-Null
$d.Content.Find.Execute("41+8=49", $true, $false, $false, $false, $false, $true, 1, $false, "50+6=56", 2) | Out-Null
$d.Content.Find.Execute("80-40=40", $true, $false, $false, $false, $false, $true, 1, $false, "18+75=93", 2) | Out-Null
$d.Content.Find.Execute("68-40=28", $true, $false, $false, $false, $false, $true, 1, $false, "35+37=72", 2) | Out-Null
$d.Content.Find.Execute("97-71=26", $true, $false, $false, $false, $false, $true, 1, $false, "97-56=41", 2) | Out-Null
$d.Content.Find.Execute("41+50=91", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=73", 2) | Out-Null
$d.Content.Find.Execute("15+68=83", $true, $false, $false, $false, $false, $true, 1, $false, "0+31=31", 2) | Out-Null
$d.Content.Find.Execute("78-13=65", $true, $false, $false, $false, $false, $true, 1, $false, "84-27=57", 2) | Out-Null
$d.Content.Find.Execute("47+7=54", $true, $false, $false, $false, $false, $true, 1, $false, "43-31=12", 2) | Out-Null
$d.Content.Find.Execute("95-11=84", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2) | Out-Null
$d.Content.Find.Execute("13+12=25", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=2", 2) | Out-Null
$d.Content.Find.Execute("61-19=42", $true, $false, $false, $false, $false, $true, 1, $false, "6+10=16", 2) | Out-Null
$d.Content.Find.Execute("83-7=76", $true, $false, $false, $false, $false, $true, 1, $false, "53-21=32", 2) | Out-Null
$d.Content.Find.Execute("3+1=4", $true, $false, $false, $false, $false, $true, 1, $false, "75-65=10", 2) | Out-Null
$d.Content.Find.Execute("58+13=71", $true, $false, $false, $false, $false, $true, 1, $false, "55-31=24", 2) | Out-Null
$d.Content.Find.Execute("93-60=33", $true, $false, $false, $false, $false, $true, 1, $false, "21+57=78", 2) | Out-Null
$d.Content.Find.Execute("4+85=89", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=90", 2) | Out-Null
$d.Content.Find.Execute("68-4=64", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=34", 2) | Out-Null
$d.Content.Find.Execute("82-59=23", $true, $false, $false, $false, $false, $true, 1, $false, "84-3=81", 2) | Out-Null
$d.Content.Find.Execute("16-11=5", $true, $false, $false, $false, $false, $true, 1, $false, "95-18=77", 2) | Out-Null
$d.Content.Find.Execute("56-16=40", $true, $false, $false, $false, $false, $true, 1, $false, "3+41=44", 2) | Out-Null
